$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 20: "2013-07-xx" / "Geodeta cz.1" / 0 / 400 / =C20+D20 ---
$ws.Range("A20").Value = "2013-07-xx"
$ws.Range("B20").Value = "Geodeta cz.1"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 400
$ws.Range("E20").Formula = "=C20+D20"

# --- Row 21: "2013-07-xx" / "Geodeta cz.2" / 0 / 650 / =C21+D21 ---
$ws.Range("A21").Value = "2013-07-xx"
$ws.Range("B21").Value = "Geodeta cz.2"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 650
$ws.Range("E21").Formula = "=C21+D21"

# --- Row 22: "2013-07-xx" / "Piasek" / 2000 / 5800 / =C22+D22 ---
$ws.Range("A22").Value = "2013-07-xx"
$ws.Range("B22").Value = "Piasek"
$ws.Range("C22").Value = 2000
$ws.Range("D22").Value = 5800
$ws.Range("E22").Formula = "=C22+D22"

# Match the formatting used by the neighbouring filled-in rows (no bold,
# bordered, numeric cells formatted as 0.00) instead of the leftover
# "empty template row" bold styling.
$ws.Range("B20:B22").Font.Bold = $false
$ws.Range("C20:C22").Font.Bold = $false

# --- View: drop the pinned topLeftCell and move the selection to B22 ---
$null = $ws.Range("B22").Select()
